$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("H2").Value = 'bedrooms'
$ws.Range("L2").Value = 'stimuli/img_72fmj.png'
$ws.Range("M2").Value = 53.87179487179487
$ws.Range("N2").Value = 36.02564102564103
$ws.Range("O2").Value = 44.94871794871795
$ws.Range("P2").Value = 39
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 3
$ws.Range("S2").Value = 3
$ws.Range("C3").Value = 1
$ws.Range("F3").Value = 2
$ws.Range("H3").Value = 'kitchens'
$ws.Range("L3").Value = 'stimuli/img_q1ynd.png'
$ws.Range("M3").Value = 70.05714285714286
$ws.Range("N3").Value = 47.31428571428572
$ws.Range("O3").Value = 58.68571428571429
$ws.Range("P3").Value = 35
$ws.Range("Q3").Value = 5
$ws.Range("R3").Value = 5
$ws.Range("S3").Value = 5
$ws.Range("C4").Value = 1
$ws.Range("F4").Value = 3
$ws.Range("H4").Value = 'kitchens'
$ws.Range("L4").Value = 'stimuli/img_411xa.png'
$ws.Range("M4").Value = 51.03030303030303
$ws.Range("N4").Value = 28.93939393939394
$ws.Range("O4").Value = 39.98484848484848
$ws.Range("P4").Value = 33
$ws.Range("Q4").Value = 2
$ws.Range("R4").Value = 2
$ws.Range("S4").Value = 2
$ws.Range("C5").Value = 1
$ws.Range("F5").Value = 4
$ws.Range("I5").Value = 'distractor'
$ws.Range("K5").Value = 'f'
$ws.Range("L5").Value = 'stimuli/img_xguy9.png'
$ws.Range("M5").Value = 78.21621621621621
$ws.Range("N5").Value = 57.24324324324324
$ws.Range("O5").Value = 67.72972972972973
$ws.Range("P5").Value = 37
$ws.Range("Q5").Value = 7
$ws.Range("R5").Value = 7
$ws.Range("S5").Value = 7
$ws.Range("C6").Value = 1
$ws.Range("F6").Value = 5
$ws.Range("H6").Value = 'bedrooms'
$ws.Range("L6").Value = 'stimuli/img_2pnl2.png'
$ws.Range("M6").Value = 6.621621621621622
$ws.Range("N6").Value = 7.135135135135135
$ws.Range("O6").Value = 6.878378378378379
$ws.Range("P6").Value = 37
$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = 1
$ws.Range("S6").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("F7").Value = 6
$ws.Range("I7").Value = 'distractor'
$ws.Range("K7").Value = 'f'
$ws.Range("L7").Value = 'stimuli/img_eppte.png'
$ws.Range("M7").Value = 78.42424242424242
$ws.Range("N7").Value = 57.03030303030303
$ws.Range("O7").Value = 67.72727272727272
$ws.Range("P7").Value = 33
$ws.Range("Q7").Value = 7
$ws.Range("R7").Value = 7
$ws.Range("S7").Value = 7
$ws.Range("C8").Value = 1
$ws.Range("F8").Value = 7
$ws.Range("H8").Value = 'bedrooms'
$ws.Range("L8").Value = 'stimuli/img_ic3os.png'
$ws.Range("M8").Value = 84.79069767441861
$ws.Range("N8").Value = 66.16279069767442
$ws.Range("O8").Value = 75.47674418604652
$ws.Range("P8").Value = 43
$ws.Range("Q8").Value = 9
$ws.Range("R8").Value = 9
$ws.Range("S8").Value = 9
$ws.Range("C9").Value = 1
$ws.Range("F9").Value = 8
$ws.Range("H9").Value = 'bedrooms'
$ws.Range("L9").Value = 'stimuli/img_anzgh.png'
$ws.Range("M9").Value = 75.10526315789474
$ws.Range("N9").Value = 55.76315789473684
$ws.Range("O9").Value = 65.4342105263158
$ws.Range("P9").Value = 38
$ws.Range("Q9").Value = 6
$ws.Range("R9").Value = 6
$ws.Range("S9").Value = 6
$ws.Range("C10").Value = 1
$ws.Range("F10").Value = 9
$ws.Range("H10").Value = 'bedrooms'
$ws.Range("L10").Value = 'stimuli/img_z3yzz.png'
$ws.Range("M10").Value = 71.71052631578948
$ws.Range("N10").Value = 49.81578947368421
$ws.Range("O10").Value = 60.76315789473685
$ws.Range("P10").Value = 38
$ws.Range("Q10").Value = 5
$ws.Range("R10").Value = 5
$ws.Range("S10").Value = 5
$ws.Range("C11").Value = 1
$ws.Range("F11").Value = 10
$ws.Range("H11").Value = 'bedrooms'
$ws.Range("L11").Value = 'stimuli/img_cmyvx.png'
$ws.Range("M11").Value = 64.25
$ws.Range("N11").Value = 40.09375
$ws.Range("O11").Value = 52.171875
$ws.Range("P11").Value = 32
$ws.Range("Q11").Value = 4
$ws.Range("R11").Value = 4
$ws.Range("S11").Value = 4
$ws.Range("C12").Value = 1
$ws.Range("F12").Value = 11
$ws.Range("H12").Value = 'bedrooms'
$ws.Range("I12").Value = 'target'
$ws.Range("K12").Value = 'j'
$ws.Range("L12").Value = 'stimuli/img_kzg3h.png'
$ws.Range("M12").Value = 77.02777777777777
$ws.Range("N12").Value = 56.22222222222222
$ws.Range("O12").Value = 66.625
$ws.Range("P12").Value = 36
$ws.Range("Q12").Value = 7
$ws.Range("R12").Value = 7
$ws.Range("S12").Value = 7
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 12
$ws.Range("H13").Value = 'bedrooms'
$ws.Range("L13").Value = 'stimuli/img_1vq1v.png'
$ws.Range("M13").Value = 69.42857142857143
$ws.Range("N13").Value = 46.59523809523809
$ws.Range("O13").Value = 58.01190476190476
$ws.Range("P13").Value = 42
$ws.Range("Q13").Value = 5
$ws.Range("R13").Value = 5
$ws.Range("S13").Value = 5
$ws.Range("C14").Value = 1
$ws.Range("F14").Value = 13
$ws.Range("H14").Value = 'bedrooms'
$ws.Range("L14").Value = 'stimuli/img_fnu4h.png'
$ws.Range("M14").Value = 85.87179487179488
$ws.Range("N14").Value = 70.71794871794872
$ws.Range("O14").Value = 78.2948717948718
$ws.Range("P14").Value = 39
$ws.Range("Q14").Value = 9
$ws.Range("C15").Value = 1
$ws.Range("F15").Value = 14
$ws.Range("H15").Value = 'bedrooms'
$ws.Range("L15").Value = 'stimuli/img_juob3.png'
$ws.Range("M15").Value = 79.92105263157895
$ws.Range("N15").Value = 59.78947368421053
$ws.Range("O15").Value = 69.85526315789474
$ws.Range("P15").Value = 38
$ws.Range("Q15").Value = 7
$ws.Range("R15").Value = 7
$ws.Range("S15").Value = 7
$ws.Range("C16").Value = 1
$ws.Range("F16").Value = 15
$ws.Range("H16").Value = 'living_rooms'
$ws.Range("I16").Value = 'distractor'
$ws.Range("K16").Value = 'f'
$ws.Range("L16").Value = 'stimuli/img_gka64.png'
$ws.Range("M16").Value = 19.23809523809524
$ws.Range("N16").Value = 20.02380952380953
$ws.Range("O16").Value = 19.63095238095238
$ws.Range("P16").Value = 42
$ws.Range("C17").Value = 1
$ws.Range("F17").Value = 16
$ws.Range("H17").Value = 'bedrooms'
$ws.Range("L17").Value = 'stimuli/img_9pfbj.png'
$ws.Range("M17").Value = 91.27272727272727
$ws.Range("N17").Value = 80.0909090909091
$ws.Range("O17").Value = 85.68181818181819
$ws.Range("P17").Value = 33
$ws.Range("Q17").Value = 10
$ws.Range("R17").Value = 10
$ws.Range("S17").Value = 10
$ws.Range("C18").Value = 1
$ws.Range("F18").Value = 17
$ws.Range("H18").Value = 'living_rooms'
$ws.Range("I18").Value = 'distractor'
$ws.Range("K18").Value = 'f'
$ws.Range("L18").Value = 'stimuli/img_16kib.png'
$ws.Range("M18").Value = 80.97727272727273
$ws.Range("N18").Value = 61.11363636363637
$ws.Range("O18").Value = 71.04545454545455
$ws.Range("P18").Value = 44
$ws.Range("Q18").Value = 8
$ws.Range("R18").Value = 8
$ws.Range("S18").Value = 8
$ws.Range("C19").Value = 1
$ws.Range("F19").Value = 18
$ws.Range("H19").Value = 'bedrooms'
$ws.Range("L19").Value = 'stimuli/img_cgdyc.png'
$ws.Range("M19").Value = 32.93023255813954
$ws.Range("N19").Value = 14.04651162790698
$ws.Range("O19").Value = 23.48837209302326
$ws.Range("P19").Value = 43
$ws.Range("Q19").Value = 1
$ws.Range("R19").Value = 1
$ws.Range("S19").Value = 1
$ws.Range("C20").Value = 1
$ws.Range("F20").Value = 19
$ws.Range("H20").Value = 'bedrooms'
$ws.Range("L20").Value = 'stimuli/img_jivhq.png'
$ws.Range("M20").Value = 37
$ws.Range("N20").Value = 22.26530612244898
$ws.Range("O20").Value = 29.63265306122449
$ws.Range("P20").Value = 49
$ws.Range("Q20").Value = 2
$ws.Range("R20").Value = 2
$ws.Range("S20").Value = 2
$ws.Range("C21").Value = 1
$ws.Range("F21").Value = 20
$ws.Range("H21").Value = 'bedrooms'
$ws.Range("I21").Value = 'target'
$ws.Range("K21").Value = 'j'
$ws.Range("L21").Value = 'stimuli/img_aweye.png'
$ws.Range("M21").Value = 53.42105263157895
$ws.Range("N21").Value = 31.84210526315789
$ws.Range("O21").Value = 42.63157894736842
$ws.Range("P21").Value = 38
$ws.Range("Q21").Value = 2
$ws.Range("R21").Value = 2
$ws.Range("S21").Value = 2
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 21
$ws.Range("H22").Value = 'bedrooms'
$ws.Range("I22").Value = 'target'
$ws.Range("K22").Value = 'j'
$ws.Range("L22").Value = 'stimuli/img_3bxjb.png'
$ws.Range("M22").Value = 87.28571428571429
$ws.Range("N22").Value = 72.65714285714286
$ws.Range("O22").Value = 79.97142857142858
$ws.Range("P22").Value = 35
$ws.Range("Q22").Value = 10
$ws.Range("R22").Value = 10
$ws.Range("S22").Value = 10
$ws.Range("C23").Value = 1
$ws.Range("F23").Value = 22
$ws.Range("H23").Value = 'bedrooms'
$ws.Range("L23").Value = 'stimuli/img_f4jxo.png'
$ws.Range("M23").Value = 82.91666666666667
$ws.Range("N23").Value = 65.52777777777777
$ws.Range("O23").Value = 74.22222222222223
$ws.Range("P23").Value = 36
$ws.Range("Q23").Value = 8
$ws.Range("R23").Value = 8
$ws.Range("S23").Value = 8
$ws.Range("C24").Value = 1
$ws.Range("F24").Value = 23
$ws.Range("H24").Value = 'bedrooms'
$ws.Range("L24").Value = 'stimuli/img_yteqw.png'
$ws.Range("M24").Value = 66.83783783783784
$ws.Range("N24").Value = 43.78378378378378
$ws.Range("O24").Value = 55.31081081081081
$ws.Range("P24").Value = 37
$ws.Range("Q24").Value = 4
$ws.Range("R24").Value = 4
$ws.Range("S24").Value = 4
$ws.Range("C25").Value = 1
$ws.Range("F25").Value = 24
$ws.Range("I25").Value = 'target'
$ws.Range("K25").Value = 'j'
$ws.Range("L25").Value = 'stimuli/img_ose78.png'
$ws.Range("M25").Value = 80.19444444444444
$ws.Range("N25").Value = 60.25
$ws.Range("O25").Value = 70.22222222222223
$ws.Range("P25").Value = 36
$ws.Range("Q25").Value = 8
$ws.Range("R25").Value = 7
$ws.Range("S25").Value = 7
$ws.Range("C26").Value = 1
$ws.Range("F26").Value = 25
$ws.Range("H26").Value = 'bedrooms'
$ws.Range("L26").Value = 'stimuli/img_t4hvr.png'
$ws.Range("M26").Value = 61.69230769230769
$ws.Range("N26").Value = 39.76923076923077
$ws.Range("O26").Value = 50.73076923076923
$ws.Range("P26").Value = 39
$ws.Range("Q26").Value = 3
$ws.Range("R26").Value = 3
$ws.Range("S26").Value = 3
$ws.Range("C27").Value = 1
$ws.Range("F27").Value = 26
$ws.Range("H27").Value = 'bedrooms'
$ws.Range("L27").Value = 'stimuli/img_gbypq.png'
$ws.Range("M27").Value = 76.275
$ws.Range("N27").Value = 51.925
$ws.Range("O27").Value = 64.1
$ws.Range("P27").Value = 40
$ws.Range("Q27").Value = 6
$ws.Range("R27").Value = 6
$ws.Range("S27").Value = 6
